# Update the cached "datetimeFigureOut" date field text (2018-12-18 -> 2018-12-19)
# on the slide master, every slide layout, and the notes master, and tweak the
# "Apache" -> "Apache Tomcat" wording on slide 6.

function Set-DatePlaceholderText($Shapes, $ShapeIndex, $NewText) {
    $sh = $Shapes.Item($ShapeIndex)
    $tr = $sh.TextFrame.TextRange
    $len = $tr.Text.Length
    $c = $tr.Characters(1, $len)
    $c.Text = $NewText
}

$p = $ppt.ActivePresentation

# ---- Slide master date placeholder ----
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes 3 "2018-12-19"

# ---- Slide layout date placeholders ----
$layouts = $master.CustomLayouts

$layoutDateShape = @{
    1 = 3
    2 = 3
    3 = 3
    4 = 4
    5 = 6
    6 = 2
    7 = 1
    8 = 4
    9 = 4
    10 = 3
    11 = 3
}

for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $shapeIdx = $layoutDateShape[$i]
    Set-DatePlaceholderText $layout.Shapes $shapeIdx "2018-12-19"
}

# ---- Notes master date placeholder ----
$notesMaster = $p.NotesMaster
Set-DatePlaceholderText $notesMaster.Shapes 2 "2018-12-19"

# ---- Slide 6: "Apache, Front-end : Bootstrap" -> "Apache Tomcat, Front-end : Bootstrap" ----
$slide6 = $p.Slides.Item(6)
$shape = $slide6.Shapes.Item(8)
$tr = $shape.TextFrame.TextRange
$fullText = $tr.Text
$needle = "Apache, Front-end : Bootstrap"
$startPos = $fullText.IndexOf($needle)
if ($startPos -ge 0) {
    $target = $tr.Characters($startPos + 1, $needle.Length)
    $target.Text = "Apache Tomcat, Front-end : Bootstrap"
}
